# Added Test Data For Hungary/Russia/Finland Market
#
# Mirrors the authoring flow: duplicate the "Croatia" country-template sheet
# (same A1:D11 layout/styles used by the other single-market sheets) three
# times, append the copies at the end of the workbook, rename them, and fill
# in each country's ticket reference (B2) + market name (B4). The last new
# sheet ("Hungary") is left as the active/selected tab with its own cursor
# position, matching the saved workbook state.

$wb = $excel.ActiveWorkbook
$template = $wb.Worksheets.Item("Croatia")

# --- Russia -----------------------------------------------------------
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count)) | Out-Null
$russia = $wb.Worksheets.Item($wb.Worksheets.Count)
$russia.Name = "Russia"
$russia.Range("B2").Value = "NGC-2929/T2908/T2922"
$russia.Range("B4").Value = "Russia Market"
$russia.Activate() | Out-Null
$russia.Range("A1:D11").Select() | Out-Null

# --- Finland ------------------------------------------------------------
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count)) | Out-Null
$finland = $wb.Worksheets.Item($wb.Worksheets.Count)
$finland.Name = "Finland"
$finland.Range("B2").Value = "NGC-3130/T2885/T2954"
$finland.Range("B4").Value = "Finland Market"
$finland.Activate() | Out-Null
$finland.Range("A1:D11").Select() | Out-Null

# --- Hungary ------------------------------------------------------------
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count)) | Out-Null
$hungary = $wb.Worksheets.Item($wb.Worksheets.Count)
$hungary.Name = "Hungary"
$hungary.Range("B2").Value = "NGC-3104/T2977"
$hungary.Range("B4").Value = "Hungary Market"
$hungary.Activate() | Out-Null
$hungary.Range("H19").Select() | Out-Null
